$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "CasesTab" Neo4j query (cell B2) had an extra `Cohort` column that
# shouldn't be part of this report. Remove the trailing
# `coalesce(co.cohort_description, '') AS `Cohort`` line (and the comma that
# used to precede it after `Response to Treatment`).
$casesQuery = @"
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
WHERE diag.disease_term IN ['Glioma']
RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,
        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,
        coalesce(s.clinical_study_type, '') AS  ``Study Type``,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,
        coalesce(demo.weight, '') AS ``Weight (kg)``,
        coalesce(diag.best_response, '') AS ``Response to Treatment``
"@
$ws.Range("B2").Value = $casesQuery

# The sheet was re-opened/re-saved at 100% zoom (previously 55%) with the
# active cell moved from B4 up to B2 (the row whose query text just changed).
$excel.ActiveWindow.Zoom = 100
$ws.Range("B2").Select()
